# Trade #34 closed at 2026-02-17 13:22:57 - unknown UNKNOWN +0.000%
#
# Adds the newly-closed trade #34 (MarketMaking strategy) to the
# "All Trades" and "MarketMaking" sheets, and rolls the aggregate
# numbers on "Summary" and "Strategy Status" forward to reflect it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1198.91   # Current Capital
$summary.Range("B4").Value = -1.09     # Total P&L $
$summary.Range("B5").Value = -0.64     # Total P&L %
$summary.Range("B6").Value = 34        # Total Trades
$summary.Range("B7").Value = 14        # Winning Trades
$summary.Range("B9").Value = 41.18     # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 98.91      # Capital
$status.Range("D4").Value = 34         # Trades
$status.Range("E4").Value = -1.09      # P&L $
$status.Range("F4").Value = -1.09      # P&L %
$status.Range("G4").Value = 41.18      # Win Rate %

# ---------------------------------------------------------------
# New trade row (#34) appended to both "All Trades" and
# "MarketMaking" sheets as row 35.
# ---------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(35, 1).Value = 34                # Trade #

    # Date column: write as literal text (leading apostrophe forces
    # text entry instead of Excel auto-parsing it into a date
    # serial), then drop back to the Normal style so no numeric
    # date format sticks to the cell.
    $ws.Cells.Item(35, 2).Value = "'2026-02-17"
    $ws.Cells.Item(35, 2).Style = "Normal"

    $ws.Cells.Item(35, 3).Value = "13:22:51"                          # Time
    $ws.Cells.Item(35, 4).Value = "MarketMaking"                      # Strategy
    $ws.Cells.Item(35, 5).Value = "UP"                                # Side
    $ws.Cells.Item(35, 6).Value = 0.69                                # Entry Price
    $ws.Cells.Item(35, 7).Value = 0.746298                            # Exit Price
    $ws.Cells.Item(35, 8).Value = "CLOSED"                            # Status
    $ws.Cells.Item(35, 9).Value = 8.1592                              # P&L %
    $ws.Cells.Item(35, 10).Value = 0.06                               # P&L $
    $ws.Cells.Item(35, 11).Value = 98.91                              # Capital After
    $ws.Cells.Item(35, 12).Value = 0                                  # Entry Slippage (bps)
    $ws.Cells.Item(35, 13).Value = 0                                  # Exit Slippage (bps)
    $ws.Cells.Item(35, 14).Value = 0.6                                # Confidence
    $ws.Cells.Item(35, 15).Value = "Normal spread capture: 19600 bps" # Entry Reason
    $ws.Cells.Item(35, 16).Value = "early_exit"                       # Exit Reason
    $ws.Cells.Item(35, 17).Value = 0.14                               # Duration (min)
}
